$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:G1").ColumnWidth = 14.1

$ws.Range("F2").Value = "vlan centos"
$ws.Range("G2").Value = "vlan window"

$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 20

$ws.Range("F4").Value = 30
$ws.Range("G4").Value = 40

$ws.Range("F5").Value = 50
$ws.Range("G5").Value = 60

$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 80

$ws.Range("F7").Value = 90
$ws.Range("G7").Value = 100

$ws.Range("F8").Value = 110
$ws.Range("G8").Value = 120

$ws.Range("F9").Value = 130
$ws.Range("G9").Value = 140

$ws.Range("G10").Select()
